$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.928.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.522.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.520.13"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.132.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.03%  "
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.858.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.520.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.47%  "
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.553"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.13%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.06%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.06%  "
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "24.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.65%  "
$ws.Range("E32").Value = "  -3.11%  "
$ws.Range("E33").Value = "  -1.41%  "
$ws.Range("E34").Value = "  -4.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "30.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.52%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "161.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.897"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.34%  "
$ws.Range("E41").Value = "  -3.20%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.740.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.26%  "
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "323.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.69%  "
